# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value  = 2598
    $ws.Range("F4").Value  = 348
    $ws.Range("F5").Value  = 1453
    $ws.Range("F13").Value = 8941
    $ws.Range("F14").Value = 387
    $ws.Range("F15").Value = 2498

    if ($sheetName -eq "展览") {
        $ws.Range("F24").Value = 2139
        $ws.Range("F26").Value = 1845
        $ws.Range("F30").Value = 680
        $ws.Range("F39").Value = 1322
        $ws.Range("F41").Value = 63
        $ws.Range("F43").Value = 282
    }
    else {
        $ws.Range("F25").Value = 2139
        $ws.Range("F27").Value = 1845
        $ws.Range("F31").Value = 680
        $ws.Range("F44").Value = 1322
        $ws.Range("F47").Value = 63
        $ws.Range("F49").Value = 282
    }
}
